# feat: add 2022-Q4 data
#
# - Duplicates the existing "2022-Q3" sheet (keeping its original figures)
#   to become the archived "2022-Q3" sheet.
# - Renames/refreshes the original sheet to "2022-Q4" with the new quarter's
#   numbers.
# - Appends a "2022-Q3" row to the "总计" (totals) summary sheet.

$wb = $excel.ActiveWorkbook
$sheetTotal = $wb.Worksheets.Item(1)
$sheetQ3 = $wb.Worksheets.Item(2)

# Step 1: Duplicate the Q3 sheet, placing the copy right after it.
# The copy keeps the old Q3 figures/formatting intact.
$sheetQ3.Copy($null, $sheetQ3)
$sheetQ3Archive = $wb.Worksheets.Item(3)

# Step 2: Rename the original sheet to "2022-Q4" first (frees up the
# "2022-Q3" name), then rename the duplicate to "2022-Q3".
$sheetQ3.Name = "2022-Q4"
$sheetQ3Archive.Name = "2022-Q3"

# Step 3: Update the (now) "2022-Q4" sheet with the new quarter's figures.
# D2:G2 are stored as text (matching the sheet's existing convention for
# this row), so force a text number-format before entering the values,
# then paste back the original (default) cell format so no stray
# number-format style lingers on the cells.
$sheetQ3.Range("D2:G2").NumberFormat = "@"
$sheetQ3.Range("D2").Value = "3.92"
$sheetQ3.Range("E2").Value = "94.38"
$sheetQ3.Range("F2").Value = "1.63"
$sheetQ3.Range("G2").Value = "0.0639"
$sheetQ3.Range("H2").Value = 10

$sheetQ3.Range("C2").Copy()
$sheetQ3.Range("D2:G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 4: Update the "总计" (totals) sheet: existing row now reflects
# 2022-Q4, and a new row is appended for the archived 2022-Q3 figures.
$sheetTotal.Range("B2").Value = "2022-Q4"

$sheetTotal.Range("A2").Copy($sheetTotal.Range("A3"))
$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("B3").Value = "2022-Q3"
$sheetTotal.Range("C3").Value = 1
$sheetTotal.Range("D3").Value = 0.06
